$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.041.33'
$ws.Range('D3').Value = '1.676.84'
$ws.Range('E3').Value = '  +0.54%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.32'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.01%  '
$ws.Range('E6').Value = '  -0.98%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  +2.34%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '21.41'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +5.38%  '
$ws.Range('E10').Value = '  +0.00%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0890'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.30%  '
$ws.Range('D12').Value = '1.914.44'
$ws.Range('E12').Value = '  +0.60%  '
$ws.Range('D13').Value = '1.683.17'
$ws.Range('E13').Value = '  +0.79%  '
$ws.Range('E14').Value = '  +0.75%  '
$ws.Range('E15').Value = '  +1.34%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.36'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.54%  '
$ws.Range('D17').Value = '27.038.87'
$ws.Range('E17').Value = '  +0.54%  '
$ws.Range('E18').Value = '  +2.64%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '236.25'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.59%  '
$ws.Range('E20').Value = '  +0.60%  '
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.46'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.36%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.24'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.06%  '
$ws.Range('E24').Value = '  -2.58%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.44'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.54%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.27'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.99%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.48'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.85%  '
$ws.Range('E28').Value = '  -0.43%  '
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0497'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.44%  '
$ws.Range('E31').Value = '  -0.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.37'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.97%  '
$ws.Range('D33').Value = '1.539.78'
$ws.Range('E33').Value = '  +6.50%  '
$ws.Range('E34').Value = '  +1.28%  '
$ws.Range('E35').Value = '  +4.94%  '
$ws.Range('E36').Value = '  -1.30%  '
$ws.Range('E37').Value = '  +0.75%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.914'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.17%  '
$ws.Range('E39').Value = '  +2.23%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.04'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.00%  '
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '67.85'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.82%  '
$ws.Range('E43').Value = '  -3.60%  '
$ws.Range('E44').Value = '  -1.48%  '
$ws.Range('D45').Value = '1.820.22'
$ws.Range('E45').Value = '  +0.75%  '
$ws.Range('E46').Value = '  -0.34%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '90.33'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.44%  '
$ws.Range('E48').Value = '  -0.05%  '
$ws.Range('E49').Value = '  +1.66%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.00'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +6.22%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0508'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.21%  '
